$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 ("ajouter regle: Deux Paires") gets its test dates and result filled in,
# mirroring the format already used by row 11 (Début / Fin date columns).
$ws.Range("H11").Copy()
$ws.Range("H13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H13").Value = 43167         # 2018-03-08

$ws.Range("I11").Copy()
$ws.Range("I13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I13").Value = 43167         # 2018-03-08

$ws.Range("J13").Value = "OK"

# Saved cursor position moves to J18
$ws.Range("J18").Select()

$wb.Save()
